# "Changing run mode in Suite.xlsx"
# The Test Suite sheet's "F Suite" row (row 7) has its Runmode cell (C7)
# flipped from "N" to "Y", and the active selection moves to C2:C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Flip the run mode for the F Suite row from "N" to "Y"
$ws.Range("C7").Value = "Y"

# Update the visible selection/active cell for the sheet
[void]$ws.Range("C2:C7").Select()
